$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N7").Value = "hhh"
$ws.Range("J9").Value = "hh"

$ws.Range("N7").Select()
